$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.017.85'
$ws.Range('E2').Value = '  +0.73%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.830.32'
$ws.Range('E3').Value = '  +0.79%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9991'
$ws.Range('E4').Value = '  +0.70%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.69'
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6184'
$ws.Range('E6').Value = '  -1.29%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9986'
$ws.Range('E7').Value = '  +0.48%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07426'
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2951'
$ws.Range('E9').Value = '  +1.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.04'
$ws.Range('E10').Value = '  +0.64%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07662'
$ws.Range('E11').Value = '  +0.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.822.06'
$ws.Range('E12').Value = '  +0.81%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.991'
$ws.Range('E13').Value = '  +0.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6725'
$ws.Range('E14').Value = '  +1.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '82.64'
$ws.Range('E15').Value = '  +0.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009054'
$ws.Range('E16').Value = '  -5.96%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.890'
$ws.Range('E17').Value = '  -1.71%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '28.980.19'
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.087.89'
$ws.Range('E19').Value = '  +1.36%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '239.01'
$ws.Range('E20').Value = '  +7.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.69'
$ws.Range('E21').Value = '  +1.58%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  +0.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.191'
$ws.Range('E23').Value = '  +1.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').Value = '  +0.74%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '158.40'
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1413'
$ws.Range('E26').Value = '  +0.86%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.479'
$ws.Range('E27').Value = '  +0.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.82'
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.491'
$ws.Range('E29').Value = '  +0.34%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05593'
$ws.Range('E30').Value = '  +3.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.105'
$ws.Range('E31').Value = '  +1.94%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.123'
$ws.Range('E32').Value = '  +0.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.213'
$ws.Range('E33').Value = '  +1.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.841'
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7413'
$ws.Range('E35').Value = '  +0.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.138'
$ws.Range('E36').Value = '  +0.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.632'
$ws.Range('E37').Value = '  +1.30%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.771'
$ws.Range('E38').Value = '  +1.49%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01774'
$ws.Range('E39').Value = '  -0.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.209.43'
$ws.Range('E40').Value = '  -1.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.389'
$ws.Range('E41').Value = '  -3.36%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8978'
$ws.Range('E42').Value = '  +0.93%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9964'
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.22'
$ws.Range('E44').Value = '  +0.14%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.985.36'
$ws.Range('E45').Value = '  +1.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '65.26'
$ws.Range('E46').Value = '  +0.85%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000123'
$ws.Range('E47').Value = '  +0.29%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5092'
$ws.Range('E48').Value = '  +0.94%  '
$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4052'
$ws.Range('E49').Value = '  +0.85%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.123'
$ws.Range('E50').Value = '  +2.07%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05812'
$ws.Range('E51').Value = '  +0.73%  '
